# Insert a new price-record row at row 756 (Kiwi / Vega Modelo de Temuco
# sheet). Inserting the row shifts every existing row from 756..829 down to
# 757..830, which matches the rest of the diff (all the other hunks are
# just the pre-existing rows sliding down by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("756:756").Insert()

$ws.Range("A756").Value = 10
$ws.Range("B756").Value = "Vega Modelo de Temuco"
$ws.Range("C756").Value = "La Araucanía"
$ws.Range("D756").Value = 45132
$ws.Range("D756").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E756").Value = 9
$ws.Range("F756").Value = "Fruta"
$ws.Range("G756").Value = 100101
$ws.Range("H756").Value = "Berries"
$ws.Range("I756").Value = 100101007
$ws.Range("J756").Value = "Kiwi"
$ws.Range("K756").Value = "Hayward"
$ws.Range("L756").Value = "Especial"
$ws.Range("M756").Value = 150
$ws.Range("N756").Value = 15000
$ws.Range("O756").Value = 15000
$ws.Range("P756").Value = 15000
$ws.Range("Q756").Value = "$/bandeja 10 kilos"
$ws.Range("R756").Value = "Región de O'Higgins"
$ws.Range("S756").Value = 1500
$ws.Range("T756").Value = 10
